# "Generate Report for Handoff" - refresh the handoff report with a new
# handoff id / content hash / timestamps, across the Overview, zh-cn and
# de-de sheets (cell values + the matching hyperlink display text).

$wb = $excel.ActiveWorkbook

$newId = "39dc5b4f-c4c5-4fdd-a7ad-d794fcc6a45f"
$newHash = "0533473f7e7b7e93bd4af3e73a29c46d29086d3a"

# ---- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("D2").Value = "2016-48-18 22:48:13"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newId.md"
    }
}

# ---- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("D2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-18 22:48:10"

foreach ($h in $wsZhCn.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newId.md"
    } elseif ($addr -eq '$D$2') {
        $h.TextToDisplay = "$newId.$newHash.zh-cn.xlf"
    }
}

# ---- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("D2").Value = "$newId.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-18 22:48:13"

foreach ($h in $wsDeDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newId.md"
    } elseif ($addr -eq '$D$2') {
        $h.TextToDisplay = "$newId.$newHash.de-de.xlf"
    }
}
